$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold + thin border + center/top alignment) from D1 to E1
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New column "name_2": Korean champion display name with internal spaces removed
$ws.Range("E1").Value = "name_2"
$ws.Range("E2").Value = "애니"
$ws.Range("E3").Value = "올라프"
$ws.Range("E4").Value = "갈리오"
$ws.Range("E5").Value = "트위스티드페이트"
$ws.Range("E6").Value = "신짜오"
$ws.Range("E7").Value = "우르곳"
$ws.Range("E8").Value = "르블랑"
$ws.Range("E9").Value = "블라디미르"
$ws.Range("E10").Value = "피들스틱"
$ws.Range("E11").Value = "케일"
$ws.Range("E12").Value = "마스터이"
$ws.Range("E13").Value = "알리스타"
$ws.Range("E14").Value = "라이즈"
$ws.Range("E15").Value = "사이온"
$ws.Range("E16").Value = "시비르"
$ws.Range("E17").Value = "소라카"
$ws.Range("E18").Value = "티모"
$ws.Range("E19").Value = "트리스타나"
$ws.Range("E20").Value = "워윅"
$ws.Range("E21").Value = "누누"
$ws.Range("E22").Value = "미스포츈"
$ws.Range("E23").Value = "애쉬"
$ws.Range("E24").Value = "트린다미어"
$ws.Range("E25").Value = "잭스"
$ws.Range("E26").Value = "모르가나"
$ws.Range("E27").Value = "질리언"
$ws.Range("E28").Value = "신지드"
$ws.Range("E29").Value = "이블린"
$ws.Range("E30").Value = "트위치"
$ws.Range("E31").Value = "카서스"
$ws.Range("E32").Value = "초가스"
$ws.Range("E33").Value = "아무무"
$ws.Range("E34").Value = "람머스"
$ws.Range("E35").Value = "애니비아"
$ws.Range("E36").Value = "샤코"
$ws.Range("E37").Value = "문도박사"
$ws.Range("E38").Value = "소나"
$ws.Range("E39").Value = "카사딘"
$ws.Range("E40").Value = "이렐리아"
$ws.Range("E41").Value = "잔나"
$ws.Range("E42").Value = "갱플랭크"
$ws.Range("E43").Value = "코르키"
$ws.Range("E44").Value = "카르마"
$ws.Range("E45").Value = "타릭"
$ws.Range("E46").Value = "베이가"
$ws.Range("E47").Value = "트런들"
$ws.Range("E48").Value = "스웨인"
$ws.Range("E49").Value = "케이틀린"
$ws.Range("E50").Value = "블리츠크랭크"
$ws.Range("E51").Value = "말파이트"
$ws.Range("E52").Value = "카타리나"
$ws.Range("E53").Value = "녹턴"
$ws.Range("E54").Value = "마오카이"
$ws.Range("E55").Value = "레넥톤"
$ws.Range("E56").Value = "자르반4세"
$ws.Range("E57").Value = "엘리스"
$ws.Range("E58").Value = "오리아나"
$ws.Range("E59").Value = "오공"
$ws.Range("E60").Value = "브랜드"
$ws.Range("E61").Value = "리신"
$ws.Range("E62").Value = "베인"
$ws.Range("E63").Value = "럼블"
$ws.Range("E64").Value = "카시오페아"
$ws.Range("E65").Value = "스카너"
$ws.Range("E66").Value = "하이머딩거"
$ws.Range("E67").Value = "나서스"
$ws.Range("E68").Value = "니달리"
$ws.Range("E69").Value = "우디르"
$ws.Range("E70").Value = "뽀삐"
$ws.Range("E71").Value = "그라가스"
$ws.Range("E72").Value = "판테온"
$ws.Range("E73").Value = "이즈리얼"
$ws.Range("E74").Value = "모데카이저"
$ws.Range("E75").Value = "요릭"
$ws.Range("E76").Value = "아칼리"
$ws.Range("E77").Value = "케넨"
$ws.Range("E78").Value = "가렌"
$ws.Range("E79").Value = "레오나"
$ws.Range("E80").Value = "말자하"
$ws.Range("E81").Value = "탈론"
$ws.Range("E82").Value = "리븐"
$ws.Range("E83").Value = "코그모"
$ws.Range("E84").Value = "쉔"
$ws.Range("E85").Value = "럭스"
$ws.Range("E86").Value = "제라스"
$ws.Range("E87").Value = "쉬바나"
$ws.Range("E88").Value = "아리"
$ws.Range("E89").Value = "그레이브즈"
$ws.Range("E90").Value = "피즈"
$ws.Range("E91").Value = "볼리베어"
$ws.Range("E92").Value = "렝가"
$ws.Range("E93").Value = "바루스"
$ws.Range("E94").Value = "노틸러스"
$ws.Range("E95").Value = "빅토르"
$ws.Range("E96").Value = "세주아니"
$ws.Range("E97").Value = "피오라"
$ws.Range("E98").Value = "직스"
$ws.Range("E99").Value = "룰루"
$ws.Range("E100").Value = "드레이븐"
$ws.Range("E101").Value = "헤카림"
$ws.Range("E102").Value = "카직스"
$ws.Range("E103").Value = "다리우스"
$ws.Range("E104").Value = "제이스"
$ws.Range("E105").Value = "리산드라"
$ws.Range("E106").Value = "다이애나"
$ws.Range("E107").Value = "퀸"
$ws.Range("E108").Value = "신드라"
$ws.Range("E109").Value = "아우렐리온솔"
$ws.Range("E110").Value = "케인"
$ws.Range("E111").Value = "조이"
$ws.Range("E112").Value = "자이라"
$ws.Range("E113").Value = "카이사"
$ws.Range("E114").Value = "나르"
$ws.Range("E115").Value = "자크"
$ws.Range("E116").Value = "야스오"
$ws.Range("E117").Value = "벨코즈"
$ws.Range("E118").Value = "탈리야"
$ws.Range("E119").Value = "카밀"
$ws.Range("E120").Value = "브라움"
$ws.Range("E121").Value = "진"
$ws.Range("E122").Value = "킨드레드"
$ws.Range("E123").Value = "징크스"
$ws.Range("E124").Value = "탐켄치"
$ws.Range("E125").Value = "루시안"
$ws.Range("E126").Value = "제드"
$ws.Range("E127").Value = "클레드"
$ws.Range("E128").Value = "에코"
$ws.Range("E129").Value = "바이"
$ws.Range("E130").Value = "아트록스"
$ws.Range("E131").Value = "나미"
$ws.Range("E132").Value = "아지르"
$ws.Range("E133").Value = "쓰레쉬"
$ws.Range("E134").Value = "일라오이"
$ws.Range("E135").Value = "렉사이"
$ws.Range("E136").Value = "아이번"
$ws.Range("E137").Value = "칼리스타"
$ws.Range("E138").Value = "바드"
$ws.Range("E139").Value = "라칸"
$ws.Range("E140").Value = "자야"
$ws.Range("E141").Value = "오른"
$ws.Range("E142").Value = "파이크"

# Page setup: A4(9)-size paper, portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the author's scroll position / active-cell selection
$ws.Range("F78").Select() | Out-Null
